$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format before assigning values so Excel keeps these
# numeric-looking / percent-looking strings as literal text (matching the
# source workbook, which stores every data cell as an inline string).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "303.23"
$ws.Range("E2").Value = "5.35%"
$ws.Range("D3").Value = "34.73"
$ws.Range("E3").Value = "12.18%"
$ws.Range("D4").Value = "5.180"
$ws.Range("E4").Value = "5.38%"
$ws.Range("D5").Value = "0.07795"
$ws.Range("E5").Value = "6.56%"
$ws.Range("D6").Value = "2.293"
$ws.Range("E6").Value = "-2.94%"
$ws.Range("D7").Value = "8.010"
$ws.Range("E7").Value = "3.63%"
$ws.Range("D8").Value = "3.991"
$ws.Range("E8").Value = "7.23%"
$ws.Range("D9").Value = "0.9283"
$ws.Range("E9").Value = "2.86%"
$ws.Range("D10").Value = "0.1012"
$ws.Range("E10").Value = "10.03%"
$ws.Range("D11").Value = "0.1815"
$ws.Range("E11").Value = "7.49%"
$ws.Range("D12").Value = "0.08504"
$ws.Range("E12").Value = "4.25%"
$ws.Range("D13").Value = "0.03472"
$ws.Range("E13").Value = "11.03%"
$ws.Range("D14").Value = "0.09902"
$ws.Range("E14").Value = "-0.32%"
$ws.Range("D15").Value = "0.001483"
$ws.Range("E15").Value = "-1.14%"
$ws.Range("D16").Value = "0.04607"
$ws.Range("E16").Value = "2.85%"
$ws.Range("D17").Value = "0.005836"
$ws.Range("E17").Value = "1.89%"
$ws.Range("D18").Value = "3.466"
$ws.Range("E18").Value = "-0.85%"
$ws.Range("E19").Value = "0.41%"
$ws.Range("D20").Value = "0.3442"
$ws.Range("E20").Value = "3.37%"
$ws.Range("D21").Value = "0.1325"
$ws.Range("E21").Value = "-0.36%"
$ws.Range("D22").Value = "4.546"
$ws.Range("E22").Value = "8.19%"
$ws.Range("D23").Value = "0.2336"
$ws.Range("E23").Value = "11.23%"
$ws.Range("D24").Value = "0.001220"
$ws.Range("E24").Value = "0.60%"
$ws.Range("D25").Value = "0.004434"
$ws.Range("E25").Value = "6.61%"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").Value = "-0.06%"
$ws.Range("D27").Value = "0.0003400"
$ws.Range("E27").Value = "0.16%"
$ws.Range("D39").Value = "0.01762"
$ws.Range("E39").Value = "11.99%"
$ws.Range("D40").Value = "0.04722"
$ws.Range("E40").Value = "6.32%"
$ws.Range("D41").Value = "0.007670"
$ws.Range("E41").Value = "4.42%"
$ws.Range("E42").Value = "6.23%"
$ws.Range("D43").Value = "0.007046"
$ws.Range("E43").Value = "-25.99%"
$ws.Range("D44").Value = "0.002299"
$ws.Range("E44").Value = "3.54%"
$ws.Range("D45").Value = "0.009818"
$ws.Range("E45").Value = "9.15%"
$ws.Range("D46").Value = "0.00005962"
$ws.Range("E46").Value = "-2.27%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("E48").Value = "9.68%"
$ws.Range("D49").Value = "0.002698"
$ws.Range("E49").Value = "34.79%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.07%"

